$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# --- Row 17: fill the previously-empty columns with the literal text "nan"
# (mirrors the "nan" placeholder already used by sibling rows), keep the
# existing Date / Correction / Serviced-by values untouched.
$ws.Range("B17").Value = "nan"
$ws.Range("C17").Value = "nan"
$ws.Range("D17").Value = "nan"
$ws.Range("E17").Value = "nan"
$ws.Range("F17").Value = "nan"
$ws.Range("G17").Value = "nan"
$ws.Range("H17").Value = "nan"
$ws.Range("I17").Value = "nan"
$ws.Range("J17").Value = "nan"
$ws.Range("K17").Value = "nan"
$ws.Range("M17").Value = "nan"

# --- Row 18: brand-new service-log entry ("Card20" event).
# A leading apostrophe forces plain-text storage: "'20" keeps "20" as text
# (matching column A elsewhere on the sheet, instead of a number), and a
# bare "'" keeps the unused columns as empty *text* cells (matching the
# sheet's existing "blank inlineStr" convention) rather than truly blank
# numeric cells.
$ws.Range("A18").Value = "'20"
$ws.Range("B18").Value = "'"
$ws.Range("C18").Value = "'"
$ws.Range("D18").Value = "'"
$ws.Range("E18").Value = "'"
$ws.Range("F18").Value = "'"
$ws.Range("G18").Value = "'"
$ws.Range("H18").Value = "'"
$ws.Range("I18").Value = "'"
$ws.Range("J18").Value = "'"
$ws.Range("K18").Value = "'"
$ws.Range("L18").Value = "14\10\2024"
$ws.Range("M18").Value = "'"
$ws.Range("N18").Value = "تم تغيير السستم من ax الي ay"
$ws.Range("O18").Value = "الخبير"

# The leading apostrophe above flips Excel's "quote-prefix" cell attribute
# on, which would otherwise persist as a oneoff style distinct from the
# plain, unstyled cells used throughout the rest of the sheet. Reset the
# style on the text-forced cells back to Normal so they stay visually/
# structurally identical to their neighbours.
$ws.Range("A18:K18").Style = "Normal"
$ws.Range("M18").Style = "Normal"
